$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.100.28'
$ws.Range("E2").Value = '  -1.53%  '

# Row 3
$ws.Range("D3").Value = '2.260.40'
$ws.Range("E3").Value = '  -1.57%  '

# Row 4
$ws.Range("E4").Value = '  +0.22%  '

# Row 5
$ws.Range("D5").Value = '''111.72'
$ws.Range("E5").Value = '  +2.76%  '

# Row 6
$ws.Range("D6").Value = '''263.76'
$ws.Range("E6").Value = '  -2.90%  '

# Row 7
$ws.Range("D7").Value = '''0.615'
$ws.Range("E7").Value = '  -1.79%  '

# Row 8
$ws.Range("E8").Value = '  +0.12%  '

# Row 9
$ws.Range("D9").Value = '''0.598'
$ws.Range("E9").Value = '  -2.99%  '

# Row 10
$ws.Range("D10").Value = '''47.40'
$ws.Range("E10").Value = '  +0.79%  '

# Row 11
$ws.Range("D11").Value = '''0.0921'
$ws.Range("E11").Value = '  -1.61%  '

# Row 12
$ws.Range("D12").Value = '''8.71'
$ws.Range("E12").Value = '  +3.88%  '

# Row 13
$ws.Range("E13").Value = '  -0.62%  '

# Row 14
$ws.Range("D14").Value = '''15.39'
$ws.Range("E14").Value = '  -2.09%  '

# Row 15
$ws.Range("D15").Value = '2.599.12'
$ws.Range("E15").Value = '  -1.50%  '

# Row 16
$ws.Range("D16").Value = '''0.849'
$ws.Range("E16").Value = '  -0.97%  '

# Row 17
$ws.Range("D17").Value = '2.257.62'
$ws.Range("E17").Value = '  -1.48%  '

# Row 18
$ws.Range("D18").Value = '43.088.44'
$ws.Range("E18").Value = '  -1.63%  '

# Row 19
$ws.Range("E19").Value = '  -3.28%  '

# Row 20
$ws.Range("D20").Value = '''6.92'
$ws.Range("E20").Value = '  +9.80%  '

# Row 21
$ws.Range("D21").Value = '''70.71'
$ws.Range("E21").Value = '  -2.08%  '

# Row 22
$ws.Range("E22").Value = '  -3.82%  '

# Row 23
$ws.Range("D23").Value = '''9.71'
$ws.Range("E23").Value = '  +4.46%  '

# Row 24
$ws.Range("D24").Value = '''229.60'
$ws.Range("E24").Value = '  -1.75%  '

# Row 25
$ws.Range("D25").Value = '''2.81'
$ws.Range("E25").Value = '  -4.65%  '

# Row 26
$ws.Range("E26").Value = '  -0.09%  '

# Row 27
$ws.Range("D27").Value = '''11.22'
$ws.Range("E27").Value = '  -1.03%  '

# Row 28
$ws.Range("D28").Value = '''3.86'
$ws.Range("E28").Value = '  -1.87%  '

# Row 29
$ws.Range("D29").Value = '''41.03'
$ws.Range("E29").Value = '  +0.69%  '

# Row 30
$ws.Range("D30").Value = '''3.38'
$ws.Range("E30").Value = '  -2.00%  '

# Row 31
$ws.Range("E31").Value = '  -1.55%  '

# Row 32
$ws.Range("D32").Value = '''171.17'
$ws.Range("E32").Value = '  -3.79%  '

# Row 33
$ws.Range("E33").Value = '  -3.32%  '

# Row 34
$ws.Range("D34").Value = '''0.0894'
$ws.Range("E34").Value = '  -2.17%  '

# Row 35
$ws.Range("E35").Value = '  -0.66%  '

# Row 36
$ws.Range("E36").Value = '  -0.89%  '

# Row 37
$ws.Range("D37").Value = '''4.60'
$ws.Range("E37").Value = '  -5.71%  '

# Row 38
$ws.Range("D38").Value = '''0.0347'
$ws.Range("E38").Value = '  -3.36%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.103'
$ws.Range("E39").Value = '  -9.01%  '

# Row 40
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = '''3.73'
$ws.Range("E40").Value = '  +2.23%  '

# Row 41
$ws.Range("D41").Value = '''14.02'
$ws.Range("E41").Value = '  +15.21%  '

# Row 42
$ws.Range("D42").Value = '''74.27'
$ws.Range("E42").Value = '  +11.59%  '

# Row 43
$ws.Range("E43").Value = '  +2.83%  '

# Row 44
$ws.Range("D44").Value = '''0.232'
$ws.Range("E44").Value = '  -1.71%  '

# Row 45
$ws.Range("D45").Value = '''6.08'
$ws.Range("E45").Value = '  +10.72%  '

# Row 46
$ws.Range("E46").Value = '  +0.13%  '

# Row 47
$ws.Range("E47").Value = '  -0.62%  '

# Row 48
$ws.Range("D48").Value = '''8.54'
$ws.Range("E48").Value = '  -2.94%  '

# Row 49
$ws.Range("D49").Value = '''0.0984'
$ws.Range("E49").Value = '  -3.47%  '

# Row 50
$ws.Range("E50").Value = '  -0.49%  '

# Row 51
$ws.Range("D51").Value = '''99.61'
$ws.Range("E51").Value = '  +0.22%  '
